# Add the walk back from Uxbridge to Ickenham
#
# The outward leg (Uxbridge -> Ickenham) was already logged; this adds the
# return leg, which in the original spreadsheet actually corresponds to a
# 5 mile *reduction* of the distance recorded against February (the
# walk total for that entry needed correcting once the return leg was
# accounted for). The edit touches the "Actual" increment cell for
# February (G3), which in turn recalculates the running "Total" (F3) via
# its formula "=F2+G3".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# February's actual-distance entry (G3) drops from 77.1 to 72.1 miles.
$ws.Range("G3").Value = 72.099999999999994

# The user's cursor ends up on G5 after making the edit.
$ws.Range("G5").Select() | Out-Null

# Restore the original active sheet (Chart1) so the workbook reopens on
# the chart, same as before the edit.
$wb.Worksheets.Item("Chart1").Activate() | Out-Null
